$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 9 (entire row), pushing the existing
# "Person FBI Identification ID" (Post Consolidation) row down to row 10.
# Excel's native row-insert behavior copies formatting from the row above,
# matching the style used by the sibling "Person State Identification ID" rows.
$ws.Rows.Item(9).Insert()

# Populate the new row with the PersonStateFingerprintIdentification mapping entry.
$ws.Range("A9").Value = "Person State Fingerprint ID"
$ws.Range("B9").Value = "An identification of a person based on a Fingerprint ID."
$ws.Range("C9").Value = "/CHcr-doc:CriminalHistoryConsolidationReport/nc:Person/CHcr-ext:PostConsolidationIdentifiers/j:PersonStateFingerprintIdentification/nc:IdentificationID"
